$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the date number-format (style index 1) that's already applied
# to the rest of column B by copying an existing formatted cell before
# writing the new values.
$ws.Range("B97").Copy()
$ws.Range("B98:B100").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Water year 2018
$ws.Range("A98").Value = 2018
$ws.Range("B98").Value = 41600
$ws.Range("C98").Value = 61900
$ws.Range("D98").Value = 14.09

# Water year 2019
$ws.Range("A99").Value = 2019
$ws.Range("B99").Value = 41969
$ws.Range("C99").Value = 27900
$ws.Range("D99").Value = 9.87

# Water year 2020
$ws.Range("A100").Value = 2020
$ws.Range("B100").Value = 42400
$ws.Range("C100").Value = 50000
$ws.Range("D100").Value = 12.71

# Move the view/selection down to the new first empty row, matching
# where Excel leaves the cursor after appending the new data.
$ws.Range("A101").Select()
